# Applies the "artfynd" row-shuffle update described by the diff:
#   - The observation records (entire rows, column B excluded) held in rows
#     3,4,5,6,7 get cyclically rotated: 3<-5, 4<-6, 5<-7, 6<-4, 7<-3.
#   - Rows 19 and 20 swap their entire record contents (column B excluded).
#   - Column B ("Taxonsorteringsordning") is a per-species sort key that the
#     source system recomputed; every row whose species (column F) is
#     "Revlummer", "Knarot" or "Kransrams" gets the new sort value,
#     regardless of which physical row it now lives in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "AY"

# Columns whose values must stay TEXT even though some of them look like
# numbers/dates ("8", "2026-01-28", "11:25", ...). Writing through Value2
# lets Excel auto-coerce number/date-looking strings, so force the cell's
# number format to Text first for exactly these columns.
$textCols = @("D","F","G","H","I","J","K","L","M","N","P","T","U","V","W", `
              "Y","Z","AA","AB","AC","AF","AH","AJ","AK","AM","AO","AT","AW","AX","AY")

function Get-RowValues($rowNum) {
    return $ws.Range("$firstCol$rowNum`:$lastCol$rowNum").Value2
}

function Set-TextFormat($rowNum) {
    foreach ($col in $textCols) {
        $ws.Range("$col$rowNum").NumberFormat = "@"
    }
}

# --- Step 1: snapshot every row we are about to touch, BEFORE any writes ---
$row3 = Get-RowValues 3
$row4 = Get-RowValues 4
$row5 = Get-RowValues 5
$row6 = Get-RowValues 6
$row7 = Get-RowValues 7
$row19 = Get-RowValues 19
$row20 = Get-RowValues 20

# --- Step 2: write the rotated/swapped content back ---
Set-TextFormat 3
$ws.Range("$firstCol" + "3:" + "$lastCol" + "3").Value2 = $row5

Set-TextFormat 4
$ws.Range("$firstCol" + "4:" + "$lastCol" + "4").Value2 = $row6

Set-TextFormat 5
$ws.Range("$firstCol" + "5:" + "$lastCol" + "5").Value2 = $row7

Set-TextFormat 6
$ws.Range("$firstCol" + "6:" + "$lastCol" + "6").Value2 = $row4

Set-TextFormat 7
$ws.Range("$firstCol" + "7:" + "$lastCol" + "7").Value2 = $row3

Set-TextFormat 19
$ws.Range("$firstCol" + "19:" + "$lastCol" + "19").Value2 = $row20

Set-TextFormat 20
$ws.Range("$firstCol" + "20:" + "$lastCol" + "20").Value2 = $row19

# --- Step 3: recompute the per-species "Taxonsorteringsordning" (col B) ---
# The source system bumped the sort key for three species; apply the new
# value to every row currently holding that species, wherever it now sits.
$sortUpdates = @{
    "Revlummer" = 97879
    "Knärot"    = 99014
    "Kransrams" = 99351
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $species = $ws.Range("F$r").Value2
    if ($null -ne $species -and $sortUpdates.ContainsKey([string]$species)) {
        $ws.Range("B$r").Value2 = $sortUpdates[[string]$species]
    }
}
